# QA Round 2: deep quality optimization - compliance, diversification, UX improvements
$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1) Text edits on the "AntonellaJourney" sheet
# -----------------------------------------------------------------------
$wsJourney = $wb.Worksheets.Item("AntonellaJourney")
$wsJourney.Range("B4").Value = "I want us to finish together 🖤"
$wsJourney.Range("B5").Value = "wait"
$wsJourney.Range("B8").Value = "fuckkkk 🖤"
$wsJourney.Range("B10").Value = "one second"
$wsJourney.Range("B11").Value = "I can't resist you anymore"
$wsJourney.Range("B20").Value = "did you see it? 🥺"
$wsJourney.Range("B22").Value = "one second 🖤"

# -----------------------------------------------------------------------
# 2) Split the old "cumcontrol" sheet into "cumcontrol1" (its own refreshed
#    content) and "cumcontrol2" (built from the old "dickpic" tab's content,
#    refreshed). The original "dickpic" wording is preserved by duplicating
#    it into a brand-new trailing "dickpic" tab before we repurpose it, and
#    "boosters" simply ends up shifted one slot later.
# -----------------------------------------------------------------------

# Rename the originals into their new roles first so the name is free for
# the later duplicate-and-rename step.
$wb.Worksheets.Item("cumcontrol").Name = "cumcontrol1"
$wb.Worksheets.Item("dickpic").Name = "cumcontrol2"

# "cumcontrol2" still holds the ORIGINAL, unedited dickpic text at this
# point -- duplicate it right after itself so the duplicate becomes the new
# "dickpic" tab (preserving the original wording), then rename it back.
$wsCumControl2 = $wb.Worksheets.Item("cumcontrol2")
$wsCumControl2.Copy($null, $wsCumControl2)
$wb.Worksheets.Item("cumcontrol2 (2)").Name = "dickpic"

# -----------------------------------------------------------------------
# 3) Refresh copy on "cumcontrol1" (was "cumcontrol")
# -----------------------------------------------------------------------
$wsCumControl1 = $wb.Worksheets.Item("cumcontrol1")

$wsCumControl1.Range("B2").Value = "if you finish before you see what I'm sending next you'll regret it 💜"

$wsCumControl1.Range("B3").Value = "wait wait wait... I have one more thing for you before you finish"
$wsCumControl1.Range("C3").Value = "DELAY. Send final PPV."

$wsCumControl1.Range("B4").Value = "I want to feel it at the same time... watch this first"
$wsCumControl1.Range("C4").Value = "SYNC variant. Send PPV."

$wsCumControl1.Range("B5").Value = "okay NOW we can go together... open this 🖤"
$wsCumControl1.Range("C5").Value = "SYNC. Send PPV."

$wsCumControl1.Range("B6").Value = "you better not be close already cutie... I have more to show you 💜"

$wsCumControl1.Range("B7").Value = "not yet... I said not yet 🖤"
$wsCumControl1.Range("C7").Value = "CONTROL. More PPVs to send. Create urgency to open next."

# -----------------------------------------------------------------------
# 4) Refresh copy on "cumcontrol2" (was "dickpic") -- new row names too
# -----------------------------------------------------------------------
$wsCumControl2 = $wb.Worksheets.Item("cumcontrol2")

$wsCumControl2.Range("A2").Value = "delay2"
$wsCumControl2.Range("B2").Value = "hold on just a little longer, I promise this next one is worth it"
$wsCumControl2.Range("C2").Value = "DELAY variant."

$wsCumControl2.Range("A3").Value = "delay1"
$wsCumControl2.Range("B3").Value = "don't you dare... not until you see what I just did 🖤"
$wsCumControl2.Range("C3").Value = "DELAY. Send PPV."

$wsCumControl2.Range("A4").Value = "sync2"
$wsCumControl2.Range("B4").Value = "let's do this together... but you have to open this first"
$wsCumControl2.Range("C4").Value = "SYNC variant."

$wsCumControl2.Range("A5").Value = "sync1"
$wsCumControl2.Range("B5").Value = "okay I'm ready now too... watch this with me 💜"
$wsCumControl2.Range("C5").Value = "SYNC. Send PPV."

$wsCumControl2.Range("A6").Value = "edge2"
$wsCumControl2.Range("B6").Value = "patience... the best part hasn't even happened yet"
$wsCumControl2.Range("C6").Value = "EDGE variant."

$wsCumControl2.Range("A7").Value = "edge1"
$wsCumControl2.Range("B7").Value = "slow down cutie... I'm not letting you off that easy 🖤"
$wsCumControl2.Range("C7").Value = "CONTROL."
